$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.759.40"
$ws.Range("E2").Value = "  -0.01%  "
$ws.Range("D3").Value = "2.804.09"
$ws.Range("E3").Value = "  +0.98%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "356.74"
$ws.Range("E5").Value = "  -0.22%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "109.20"
$ws.Range("E6").Value = "  -0.27%  "
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.632"
$ws.Range("E9").Value = "  +6.99%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.93"
$ws.Range("E10").Value = "  +0.20%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0840"
$ws.Range("E12").Value = "  -0.44%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "19.93"
$ws.Range("E13").Value = "  +2.29%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.77"
$ws.Range("E14").Value = "  +2.24%  "
$ws.Range("D15").Value = "3.241.02"
$ws.Range("E15").Value = "  +0.71%  "
$ws.Range("D16").Value = "2.799.66"
$ws.Range("E16").Value = "  +0.62%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.943"
$ws.Range("E17").Value = "  +0.44%  "
$ws.Range("D18").Value = "51.719.03"
$ws.Range("E18").Value = "  +0.06%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.70"
$ws.Range("E19").Value = "  +3.41%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.16"
$ws.Range("E20").Value = "  +5.12%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.59"
$ws.Range("E21").Value = "  +3.36%  "
$ws.Range("D22").Value = "0.0₃0979"
$ws.Range("E22").Value = "  +0.75%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.42"
$ws.Range("E23").Value = "  +0.34%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "268.30"
$ws.Range("E24").Value = "  -0.31%  "
$ws.Range("E25").Value = "  +0.57%  "
$ws.Range("E26").Value = "  +0.10%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "26.14"
$ws.Range("E27").Value = "  -0.93%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.164"
$ws.Range("E28").Value = "  +0.29%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.37"
$ws.Range("E29").Value = "  +1.33%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "37.94"
$ws.Range("E30").Value = "  +9.24%  "
$ws.Range("E31").Value = "  +0.65%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.20"
$ws.Range("E32").Value = "  -0.80%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "51.97"
$ws.Range("E33").Value = "  +0.26%  "
$ws.Range("E34").Value = "  +11.04%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0444"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0862"
$ws.Range("E36").Value = "  +2.45%  "
$ws.Range("E37").Value = "  -0.13%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.89"
$ws.Range("E38").Value = "  +0.73%  "
$ws.Range("E39").Value = "  +2.01%  "
$ws.Range("E40").Value = "  +0.03%  "
$ws.Range("E41").Value = "  +1.09%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.49"
$ws.Range("E42").Value = "  -2.20%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "22.04"
$ws.Range("E43").Value = "  +1.83%  "
$ws.Range("E44").Value = "  -1.28%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "119.00"
$ws.Range("E45").Value = "  -0.58%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.47"
$ws.Range("E46").Value = "  +8.36%  "
$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.39"
$ws.Range("E47").Value = "  +3.68%  "
$ws.Range("B48").Value = "Maker"
$ws.Range("C48").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D48").Value = "2.107.40"
$ws.Range("E48").Value = "  +1.18%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.38"
$ws.Range("E49").Value = "  +9.81%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.912"
$ws.Range("E50").Value = "  -2.60%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.42"
$ws.Range("E51").Value = "  -5.74%  "
